$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Vinicius-Des. Maq. C"
$ws.Range("C2").Value = "[André B.-Elet. Digit. Básica, Carlos-Tornearia, Victor-Ajustagem, Elaine-Metalografia]"
$ws.Range("E2").Value = "-"

# Row 3
$ws.Range("B3").Value = "Maria Celeste-Maq. Term. FL"
$ws.Range("C3").Value = "[André B.-Elet. Digit. Básica, Carlos-Tornearia, Victor-Ajustagem, Elaine-Metalografia]"
$ws.Range("D3").Value = "Vinicius-Des. Maq. C"
$ws.Range("E3").Value = "-"

# Row 4
$ws.Range("B4").Value = "Maria Celeste-Maq. Term. FL"
$ws.Range("C4").Value = "[André B.-Elet. Digit. Básica, Carlos-Tornearia, Victor-Ajustagem, Elaine-Metalografia]"
$ws.Range("D4").Value = "Vinicius-Des. Maq. C"
$ws.Range("E4").Value = "-"

# Row 6
$ws.Range("B6").Value = "Nilton-Mec. Tec. Res. "
$ws.Range("C6").Value = "[André B.-Elet. Digit. Básica, Carlos-Tornearia, Victor-Ajustagem, Elaine-Metalografia]"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"

# Row 7
$ws.Range("B7").Value = "Vinicius-Des. Maq. C"
$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "-"

# Row 8
$ws.Range("B8").Value = "Vinicius-Des. Maq. C"
$ws.Range("C8").Value = "-"
